# Week 1 syllabus update:
# - Remove the "Quizzes" column (F) entirely (shift cells left)
# - Update the Week 0 "Required Readings" link (D2) to the new URL
# - Move the active cell selection to D14 (matches the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "Quizzes" column (F), shifting remaining cells left.
$ws.Columns.Item(6).Delete()

# Update the reading link in D2 with the new URL.
$ws.Range("D2").Value = "LSWR Ch 2 and 3 <br> [Clayton 2020](https://nautil.us/how-eugenics-shaped-statistics-238014/)"

# Update the selected cell to match the saved view.
$ws.Range("D14").Select()
